$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Crouching
$ws.Range("D8").Value = "BossCrouching_Yellow.png"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "49x67"
$ws.Range("H8").Value = "DONE"
$ws.Rows.Item(8).RowHeight = 15

# Row 13: CrouchPunch
$ws.Range("D13").Value = "BossCrouchPunch_Yellow.png"
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = "147x67"
$ws.Range("H13").Value = "DONE"
$ws.Rows.Item(13).RowHeight = 15

# Column width tweaks
$ws.Columns.Item(4).ColumnWidth = 25.6

# Update selection
$ws.Range("G18").Select()
